# mlk holiday jan added
#
# The "16-End" sheet of the January 2023 CBOC sign-in workbook lists the
# days of the month across pairs of columns (one pair per day). Weekend
# days (SAT/SUN) are already formatted with a gray fill, narrower columns,
# and "X" marks down the column (since the clinics are closed). Monday,
# January 16 2023 is MLK Day - a holiday - so its column (B/C) needs to be
# re-formatted the same way the weekend columns (e.g. N/O, the adjacent
# SUN column for that week) already are: narrower width, gray-shaded
# borders/fill copied in, and "X" placed in every data row that the
# weekend column has an "X" in.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("16-End")

# Column widths for B (MON) and C (its paired column) should match the
# already-narrow weekend columns (e.g. N/O).
$ws.Range("B1").ColumnWidth = $ws.Range("N1").ColumnWidth
$ws.Range("C1").ColumnWidth = $ws.Range("O1").ColumnWidth

# Header rows: copy formatting only (values/text such as "MON", the date
# 16, "Tech", "Time of Arrival" must stay put - only the look changes).
$ws.Range("N2:O3").Copy()
$ws.Range("B2").PasteSpecial(-4122)

$ws.Range("N4:O4").Copy()
$ws.Range("B4").PasteSpecial(-4122)

# Data rows: copy the gray/weekend formatting down from N/O into B/C, then
# stamp an "X" into the rows that are marked off on the weekend column.
$xRows = @(5, 6, 8, 9, 11, 12, 14, 15, 17, 18, 20, 21, 23, 24, 26, 27)

for ($r = 5; $r -le 27; $r++) {
    $ws.Range("N" + $r + ":O" + $r).Copy()
    $ws.Range("B" + $r).PasteSpecial(-4122)

    if ($xRows -contains $r) {
        $ws.Range("B" + $r).Value = "X"
        $ws.Range("C" + $r).Value = "X"
    }
}

$excel.CutCopyMode = 0
